# Auto-applies the diff to Marilith_Profits workbook (per-class leve-profit sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 461.33334
$ws.Range("J29").Value = 1000
$ws.Range("L29").Value = 3000
$ws.Range("N29").Value = -3562

$ws.Range("H31").Value = 15.333333
$ws.Range("I31").Value = 15.333333
$ws.Range("K31").Value = 45.999999
$ws.Range("M31").Value = 184.000001

$ws.Range("H33").Value = 1021.8889
$ws.Range("I33").Value = 532.3333
$ws.Range("J33").Value = 2001
$ws.Range("K33").Value = 532.3333
$ws.Range("L33").Value = 2001
$ws.Range("M33").Value = -303.3333
$ws.Range("N33").Value = -2459

$ws.Range("H41").Value = 751.9
$ws.Range("I41").Value = 427.8
$ws.Range("J41").Value = 1076
$ws.Range("K41").Value = 427.8
$ws.Range("L41").Value = 1076
$ws.Range("M41").Value = 12.19999999999999
$ws.Range("N41").Value = -1956

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H48").Value = 679.5833
$ws.Range("I48").Value = 2700
$ws.Range("J48").Value = 495.9091
$ws.Range("K48").Value = 8100
$ws.Range("L48").Value = 1487.7273
$ws.Range("M48").Value = -7808
$ws.Range("N48").Value = -2071.7273

$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

$ws.Range("H54").Value = 7666.6665
$ws.Range("J54").Value = 7666.6665
$ws.Range("L54").Value = 7666.6665
$ws.Range("N54").Value = -8638.666499999999

$ws.Range("H56").Value = 679.5833
$ws.Range("I56").Value = 2700
$ws.Range("J56").Value = 495.9091
$ws.Range("K56").Value = 8100
$ws.Range("L56").Value = 1487.7273
$ws.Range("M56").Value = -7566
$ws.Range("N56").Value = -2555.7273

$ws.Range("H59").Value = 2000
$ws.Range("J59").Value = 2000
$ws.Range("L59").Value = 6000
$ws.Range("N59").Value = -7114

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

$ws.Range("H76").Value = 3596.2856
$ws.Range("I76").Value = 3135.8
$ws.Range("K76").Value = 3135.8
$ws.Range("M76").Value = -2820.8

$ws.Range("H79").Value = 3596.2856
$ws.Range("I79").Value = 3135.8
$ws.Range("K79").Value = 3135.8
$ws.Range("M79").Value = -2043.8

$ws.Range("H86").Value = 14499.143
$ws.Range("I86").Value = 17399
$ws.Range("J86").Value = 7249.5
$ws.Range("K86").Value = 17399
$ws.Range("L86").Value = 7249.5
$ws.Range("M86").Value = -16276
$ws.Range("N86").Value = -9495.5

$ws.Range("H87").Value = 37499
$ws.Range("J87").Value = 37499
$ws.Range("L87").Value = 37499
$ws.Range("N87").Value = -39995

$ws.Range("H89").Value = 14499.143
$ws.Range("I89").Value = 17399
$ws.Range("J89").Value = 7249.5
$ws.Range("K89").Value = 86995
$ws.Range("L89").Value = 36247.5
$ws.Range("M89").Value = -81379
$ws.Range("N89").Value = -47479.5

$ws.Range("H90").Value = 37499
$ws.Range("J90").Value = 37499
$ws.Range("L90").Value = 112497
$ws.Range("N90").Value = -124977

$ws.Range("H107").Value = 3329.375
$ws.Range("I107").Value = 1105
$ws.Range("J107").Value = 10002.5
$ws.Range("K107").Value = 1105
$ws.Range("L107").Value = 10002.5
$ws.Range("M107").Value = 815
$ws.Range("N107").Value = -13842.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 444
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 444
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 444
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -784

$ws.Range("H32").Value = 6493.343
$ws.Range("I32").Value = 5213.7354
$ws.Range("K32").Value = 5213.7354
$ws.Range("M32").Value = -4926.7354

$ws.Range("H61").Value = 1660
$ws.Range("I61").Value = 1660
$ws.Range("K61").Value = 1660
$ws.Range("M61").Value = -1448

$ws.Range("H88").Value = 4890.5625
$ws.Range("I88").Value = 1866
$ws.Range("J88").Value = 5898.75
$ws.Range("K88").Value = 1866
$ws.Range("L88").Value = 5898.75
$ws.Range("M88").Value = -1460
$ws.Range("N88").Value = -6710.75

$ws.Range("H91").Value = 4890.5625
$ws.Range("I91").Value = 1866
$ws.Range("J91").Value = 5898.75
$ws.Range("K91").Value = 1866
$ws.Range("L91").Value = 5898.75
$ws.Range("M91").Value = -462
$ws.Range("N91").Value = -8706.75

$ws.Range("H132").Value = 1313.5294
$ws.Range("I132").Value = 1274.3077
$ws.Range("K132").Value = 3822.9231
$ws.Range("M132").Value = -1292.9231

$ws.Range("H136").Value = 1660
$ws.Range("I136").Value = 1660
$ws.Range("K136").Value = 4980
$ws.Range("M136").Value = -2430

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1480.6
$ws.Range("I20").Value = 745
$ws.Range("K20").Value = 745
$ws.Range("M20").Value = -498

$ws.Range("H107").Value = 2299.5
$ws.Range("I107").Value = 1074.25
$ws.Range("K107").Value = 1074.25
$ws.Range("M107").Value = 845.75

$ws.Range("H132").Value = 97000
$ws.Range("J132").Value = 97000
$ws.Range("L132").Value = 97000
$ws.Range("N132").Value = -107120

$ws.Range("H134").Value = 7910.1
$ws.Range("I134").Value = 7910.1
$ws.Range("K134").Value = 23730.3
$ws.Range("M134").Value = -21195.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 2124.7144
$ws.Range("J3").Value = 4400.3335
$ws.Range("L3").Value = 4400.3335
$ws.Range("N3").Value = -4626.3335

$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 50
$ws.Range("N22").ClearContents()

$ws.Range("H99").Value = 3759.2727
$ws.Range("I99").Value = 6855
$ws.Range("J99").Value = 1179.5
$ws.Range("K99").Value = 6855
$ws.Range("L99").Value = 1179.5
$ws.Range("M99").Value = -5357
$ws.Range("N99").Value = -4175.5

$ws.Range("H126").Value = 3759.2727
$ws.Range("I126").Value = 6855
$ws.Range("J126").Value = 1179.5
$ws.Range("K126").Value = 20565
$ws.Range("L126").Value = 3538.5
$ws.Range("M126").Value = -18095
$ws.Range("N126").Value = -8478.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 468.625
$ws.Range("I2").Value = 150
$ws.Range("J2").Value = 514.1429000000001
$ws.Range("K2").Value = 900
$ws.Range("L2").Value = 3084.8574
$ws.Range("M2").Value = -787
$ws.Range("N2").Value = -3310.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2989.1428
$ws.Range("I70").Value = 2975
$ws.Range("J70").Value = 2999.75
$ws.Range("K70").Value = 2975
$ws.Range("L70").Value = 2999.75
$ws.Range("M70").Value = -2705
$ws.Range("N70").Value = -3539.75

$ws.Range("H73").Value = 2989.1428
$ws.Range("I73").Value = 2975
$ws.Range("J73").Value = 2999.75
$ws.Range("K73").Value = 2975
$ws.Range("L73").Value = 2999.75
$ws.Range("M73").Value = -2039
$ws.Range("N73").Value = -4871.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 93
$ws.Range("I58").Value = 93
$ws.Range("K58").Value = 93
$ws.Range("M58").Value = 167

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 41250
$ws.Range("I4").Value = 40000
$ws.Range("K4").Value = 40000
$ws.Range("M4").Value = -39887

$ws.Range("H13").Value = 497.5
$ws.Range("I13").Value = 497.5
$ws.Range("K13").Value = 497.5
$ws.Range("M13").Value = -357.5

$ws.Range("H32").Value = 4000
$ws.Range("I32").Value = 4000
$ws.Range("K32").Value = 4000
$ws.Range("M32").Value = -3683

$ws.Range("H51").Value = 28999
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H107").Value = 1388.1428
$ws.Range("I107").Value = 1408
$ws.Range("J107").Value = 1269
$ws.Range("K107").Value = 4224
$ws.Range("L107").Value = 3807
$ws.Range("M107").Value = -2304
$ws.Range("N107").Value = -7647

$ws.Range("H126").Value = 3117.35
$ws.Range("I126").Value = 3033.6924
$ws.Range("K126").Value = 9101.0772
$ws.Range("M126").Value = -6631.0772
